$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Heure " header to remove the trailing space
$ws.Range("C2").Value = "Heure"

# Add row 11
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Côme"
$ws.Range("C11").Value = "09:29"
$ws.Range("D11").Value = "Côme"

# Add row 12
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Baptiste"
$ws.Range("C12").Value = "09:30"
$ws.Range("D12").Value = "Mathieu"
